# This script reproduces the edit described by the diff:
#  - Sheet "展览" (index 1): the "10.19剑网3同人only" row (row 2) is removed,
#    remaining rows shift up, and several counter/price values are refreshed.
#  - Sheet "全部类型" (index 4): same "10.19剑网3同人only" row (row 2) is removed,
#    remaining rows shift up, and several counter/price values are refreshed.
#  - Sheets "演出" and "本地生活" are left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "展览" (exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Remove row 2 (南宁·10.19剑网3同人only——寒光烈火·阵营PK战); everything below
# shifts up by one row automatically.
$ws1.Rows.Item(2).Delete()

# New row 2 (was row 3): 南宁·熊喵M动漫嘉年华·万圣派对 -- refresh counts
$ws1.Cells.Item(2, 6).Value = 315
$ws1.Cells.Item(2, 7).Value = 60

# New row 3 (was row 4): 南宁·万圣漫控嘉年华10 -- refresh counts
$ws1.Cells.Item(3, 6).Value = 1318
$ws1.Cells.Item(3, 7).Value = 50

# New row 4 (was row 5): 南宁·梦中礼Lolita茶会 -- refresh counts
$ws1.Cells.Item(4, 6).Value = 86
$ws1.Cells.Item(4, 7).Value = 138

# New row 5 (was row 6): 南宁·黑塔利亚同人ONLY -- refresh counts
$ws1.Cells.Item(5, 6).Value = 66
$ws1.Cells.Item(5, 7).Value = 58

# ---------------------------------------------------------------------------
# Sheet 4: "全部类型" (all types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Remove row 2 (南宁·10.19剑网3同人only——寒光烈火·阵营PK战); everything below
# shifts up by one row automatically.
$ws4.Rows.Item(2).Delete()

# New row 3 (was row 4): 南宁·熊喵M动漫嘉年华·万圣派对 -- refresh counts
$ws4.Cells.Item(3, 6).Value = 315
$ws4.Cells.Item(3, 7).Value = 60

# New row 4 (was row 5): 南宁·万圣漫控嘉年华10 -- refresh counts
$ws4.Cells.Item(4, 6).Value = 1318
$ws4.Cells.Item(4, 7).Value = 50

# New row 5 (was row 6): 南宁·梦中礼Lolita茶会 -- refresh counts
$ws4.Cells.Item(5, 6).Value = 86
$ws4.Cells.Item(5, 7).Value = 138

# New row 6 (was row 7): 南宁·黑塔利亚同人ONLY -- refresh counts
$ws4.Cells.Item(6, 6).Value = 66
$ws4.Cells.Item(6, 7).Value = 58
